# "break out stock.yaml completed"
# 1) week sheet: rows 83-93 had their `bsecode` (column D) values written in as
#    text; convert them to real numbers.
# 2) month sheet: append the newly scraped batch of rows (35-47).

$wb = $excel.ActiveWorkbook

# --- 1. Fix column D (bsecode) on the "week" sheet: text -> number ---------
$weekSheet = $wb.Worksheets.Item("week")
$weekBsecodes = @{
    83 = 539523
    84 = 500257
    85 = 524715
    86 = 532321
    87 = 542830
    88 = 500295
    89 = 540222
    90 = 513599
    91 = 500103
    92 = 526371
    93 = 500113
}
foreach ($row in $weekBsecodes.Keys) {
    $weekSheet.Cells.Item($row, 4).Value = $weekBsecodes[$row]
}

# --- 2. Append new rows 35-47 to the "month" sheet --------------------------
$monthSheet = $wb.Worksheets.Item("month")

# sr, nsecode, name, bsecode, per_chg, close, volume, timeframe, Date Time
$newMonthRows = @(
    @(1,  "POLYCAB",    "Polycab India Ltd",                                       "542652", -3.57, 6739.5,  8853879,  "month", "30/06/2024 21:34:13"),
    @(2,  "LTIM",       "LTI Mindtree Ltd",                                        "540005",  0.15, 5385.05, 413940,   "month", "30/06/2024 21:34:13"),
    @(3,  "PIIND",      "Pi Industries Limited",                                   "523642", -0.54, 3798.75, 172804,   "month", "30/06/2024 21:34:13"),
    @(4,  "TITAN",      "Titan Company Limited",                                   "500114",  0.7,  3404.2,  1005498,  "month", "30/06/2024 21:34:13"),
    @(5,  "NESTLEIND",  "Nestle India Limited",                                    "500790",  0.71, 2551.65, 1182033,  "month", "30/06/2024 21:34:13"),
    @(6,  "ASTRAL",     "Astral Poly Technik Limited",                             "532830", -1.53, 2380.8,  333399,   "month", "30/06/2024 21:34:13"),
    @(7,  "INDUSINDBK", "Indusind Bank Limited",                                   "532187", -2.55, 1464.5,  6425323,  "month", "30/06/2024 21:34:13"),
    @(8,  "GNFC",       "Gujarat Narmada Valley Fertilizers And Chemicals Limited","500670",  5.68, 710.5,   7334980,  "month", "30/06/2024 21:34:13"),
    @(9,  "AUBANK",     "AU Small Finance Bank",                                   "540611",  0.89, 672.05,  1933315,  "month", "30/06/2024 21:34:13"),
    @(10, "ICICIPRULI", "Icici Prudential Life Insurance Company Limited",         "540133",  1.01, 605.7,   1685280,  "month", "30/06/2024 21:34:13"),
    @(11, "ITC",        "Itc Limited",                                             "500875", -0.16, 424.9,   17866326, "month", "30/06/2024 21:34:13"),
    @(12, "RBLBANK",    "Rbl Bank Limited",                                        "540065",  0.01, 262.98,  6635809,  "month", "30/06/2024 21:34:13"),
    @(13, "BANDHANBNK", "Bandhan Bank Ltd",                                        "541153",  1.74, 203.78,  9047376,  "month", "30/06/2024 21:34:13")
)

$r = 35
foreach ($rowData in $newMonthRows) {
    $monthSheet.Cells.Item($r, 1).Value = $rowData[0]
    $monthSheet.Cells.Item($r, 2).Value = $rowData[1]
    $monthSheet.Cells.Item($r, 3).Value = $rowData[2]

    # bsecode arrives (just like the old "week" rows before they were fixed)
    # as text holding a numeric-looking string - force it to stay text.
    $monthSheet.Cells.Item($r, 4).NumberFormat = "@"
    $monthSheet.Cells.Item($r, 4).Value = $rowData[3]

    $monthSheet.Cells.Item($r, 5).Value = $rowData[4]
    $monthSheet.Cells.Item($r, 6).Value = $rowData[5]
    $monthSheet.Cells.Item($r, 7).Value = $rowData[6]
    $monthSheet.Cells.Item($r, 8).Value = $rowData[7]
    $monthSheet.Cells.Item($r, 9).Value = $rowData[8]
    $r = $r + 1
}
